$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column L width: 47.71 -> 34 characters ---
$ws.Range("L1").EntireColumn.ColumnWidth = 33.166666667

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 38.25
$ws.Rows.Item(8).RowHeight = 75

# --- Date value updates (2do Trimestre 2023 report) ---
$ws.Range("B8").Value = 45017
$ws.Range("C8").Value = 45107
$ws.Range("J8").Value = 45117
$ws.Range("K8").Value = 45117

# --- Border tweaks ---
# A3 (merged A3:C3 label cell): left-only -> full box
$ws.Range("A3").Borders.Item(8).LineStyle = 1
$ws.Range("A3").Borders.Item(8).Weight = 2
$ws.Range("A3").Borders.Item(9).LineStyle = 1
$ws.Range("A3").Borders.Item(9).Weight = 2
$ws.Range("A3").Borders.Item(10).LineStyle = 1
$ws.Range("A3").Borders.Item(10).Weight = 2

# C3 (end of merged A3:C3 range): right-only -> none
$ws.Range("C3").Borders.Item(10).LineStyle = 0

# L8 (note cell): full box -> right+bottom only
$ws.Range("L8").Borders.Item(7).LineStyle = 0
$ws.Range("L8").Borders.Item(8).LineStyle = 0
$ws.Range("L8").Borders.Item(9).LineStyle = 1
$ws.Range("L8").Borders.Item(9).Weight = 2
$ws.Range("L8").Borders.Item(10).LineStyle = 1
$ws.Range("L8").Borders.Item(10).Weight = 2

# --- View state: scroll + selection ---
$win = $excel.ActiveWindow
$ws.Range("B13").Select()
$win.ScrollRow = 3
$win.ScrollColumn = 1
